$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-43: corrected open/close/high/low (D:G), shares_outstanding (H),
# and fixed_ticker (I) now collapsed to "QCOM" for every row.
$data = @(
    @(2, 50.63421817677034, 50.28664779663086, 53.17073607807491, 48.94074068422496, 1079000000),
    @(3, 46.98672627648549, 47.94731140136719, 49.18341503964993, 45.29639801189692, 1079000000),
    @(4, 40.24893760821039, 44.6192512512207, 45.94836766948758, 39.79839093213948, 1079000000),
    @(5, 37.48408630728156, 34.38253021240234, 38.37132709371375, 32.96446019572154, 1079000000),
    @(6, 39.09888826657534, 38.67023849487305, 40.78286304697232, 38.06553692049063, 1079000000),
    @(7, 40.95025544383655, 48.36133193969727, 48.54680039631238, 39.53604445929848, 1079000000),
    @(8, 53.05826797356231, 53.55705261230469, 55.81717385122092, 50.11231934287078, 1079000000),
    @(9, 51.72956938077056, 41.96645736694336, 52.68781622568582, 41.49519009685886, 1079000000),
    @(10, 45.39519688131731, 42.60457611083984, 45.60925074610093, 40.47196706406418, 1079000000),
    @(11, 44.45720446320183, 42.58380508422852, 45.5860474829653, 41.81523176312123, 1079000000),
    @(12, 42.16508008745418, 41.29085540771485, 44.33444900991759, 40.4894858724879, 1079000000),
    @(13, 52.55302815348841, 55.71208953857422, 56.24268211307936, 52.38160708686337, 1079000000),
    @(14, 45.36386671677448, 41.99656677246094, 46.75524689147601, 39.97948259651921, 1079000000),
    @(15, 46.0874189697698, 53.316650390625, 53.79084110621253, 45.85448401477702, 1079000000),
    @(16, 60.8285512684193, 52.79475021362305, 62.04579773411955, 51.35924416486456, 1079000000),
    @(17, 47.69083390378788, 42.02224349975586, 49.3625566405369, 41.66583345403984, 1079000000),
    @(18, 49.48156114839674, 73.93904876708984, 76.08519561269873, 47.89341379255026, 1079000000),
    @(19, 69.88485925290979, 63.37891006469727, 69.96283018049992, 63.03238656188953, 1079000000),
    @(20, 67.2875267411555, 70.2388916015625, 72.58775097842894, 63.57649640492343, 1079000000),
    @(21, 78.35992709617932, 75.06889343261719, 84.62519693545927, 74.83130874376995, 1079000000),
    @(22, 57.65199971216418, 69.77665710449219, 71.72795961644098, 57.25287319015598, 1079000000),
    @(23, 81.34645581221291, 94.39620971679688, 96.30898151881948, 78.21808982320253, 1079000000),
    @(24, 107.7698315357262, 110.85205078125, 118.9934199809431, 103.7261040223708, 1079000000),
    @(25, 138.0890951352035, 141.0402221679688, 151.5631906849169, 132.7915175698095, 1079000000),
    @(26, 121.9626364152752, 125.86181640625, 131.7559215118209, 120.1218631584511, 1079000000),
    @(27, 130.2499556137421, 136.5295867919922, 138.0972173070809, 123.5875328781401, 1079000000),
    @(28, 118.8523103021039, 121.8190307617188, 122.542407791879, 111.8658412218451, 1079000000),
    @(29, 168.1160033275405, 161.5442962646484, 177.9229980886852, 144.6509042416479, 1079000000),
    @(30, 140.2334316150178, 128.9272613525391, 142.254700225449, 120.5929998285631, 1079000000),
    @(31, 116.6722822082895, 134.5883178710938, 145.3509356502369, 111.6713793860902, 1079000000),
    @(32, 107.2876919650483, 109.7788543701172, 118.689168190246, 97.64962555494668, 1079000000),
    @(33, 104.5539255832556, 125.0797424316406, 126.6571993836932, 100.5914855326328, 1079000000),
    @(34, 119.9677190184511, 110.3412094116211, 120.5156479724526, 105.9483431554105, 1079000000),
    @(35, 113.4570753205749, 125.7283630371094, 126.4893759410638, 108.0824479950405, 1079000000),
    @(36, 105.7630125407569, 104.4122314453125, 109.5950056957919, 99.94796322034904, 1079000000),
    @(37, 137.0841068397152, 143.1771545410156, 152.307097833942, 130.0944468247821, 1079000000),
    @(38, 163.0516296963782, 160.7066802978516, 170.5322189180645, 151.4915941188385, 1079000000),
    @(39, 194.076135630825, 176.0569305419922, 205.3819146181004, 162.2603788552847, 1079000000),
    @(40, 166.0622141262047, 159.1869659423828, 175.8714236376172, 158.5121524493731, 1079000000),
    @(41, 152.3907915403076, 170.0079956054688, 173.802773808184, 150.2083016671886, 1079000000),
    @(42, 150.6228848329906, 146.7384643554688, 153.657292349785, 119.3992061182842, 1079000000),
    @(43, 157.3246773091798, 145.9295196533203, 163.0322016739741, 144.7562042348384, 1079000000)
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 4).Value = $item[1]   # D: open_price
    $ws.Cells.Item($r, 5).Value = $item[2]   # E: close_price
    $ws.Cells.Item($r, 6).Value = $item[3]   # F: high_price
    $ws.Cells.Item($r, 7).Value = $item[4]   # G: low_price
    $ws.Cells.Item($r, 8).Value = $item[5]   # H: shares_outstanding
    $ws.Cells.Item($r, 9).Value = "QCOM"     # I: fixed_ticker
}
